# Update "想去人数" (number of people interested) counts for three events.
# These events appear both on their originating sheet ("展览") and on the
# aggregated "全部类型" sheet, so both copies must be updated in tandem.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# 展览 sheet: rows 5, 8, 14 -> column F
$ws1.Range("F5").Value = 1047
$ws1.Range("F8").Value = 202
$ws1.Range("F14").Value = 12396

# 全部类型 sheet: rows 7, 10, 16 -> column F (same events, mirrored rows)
$ws4.Range("F7").Value = 1047
$ws4.Range("F10").Value = 202
$ws4.Range("F16").Value = 12396
